$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("AttackRangeRate"), shifting
# AttackRangeRate/AttackSpeedRate/MovementSpeed one column to the right
# (D,E,F -> E,F,G), carrying their column widths/styles along.
[void]$ws.Columns("D:D").Insert()

# Header for the new column.
$ws.Range("D2").Value = "AttackRangeForward"

# New stat column is 0 for every character row.
for ($r = 3; $r -le 12; $r++) {
    $ws.Cells.Item($r, 4).Value = 0
}

# Match the author's last on-screen selection.
[void]$ws.Range("H7").Select()
